$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update genders that changed from male to female
$ws.Range("H4").Value = "female"
$ws.Range("H7").Value = "female"
$ws.Range("H9").Value = "female"
$ws.Range("H11").Value = "female"

# 2. Insert a new row above the "freedom" row (currently row 13), pushing it to row 14
$ws.Rows("13:13").Insert()

# 3. Fill in the new row 13 with the "qinshisan" user record
$ws.Range("A13").Value = "qinshisan"
$ws.Range("B13").Value = "秦十三"
$ws.Range("C13").Value = "qinshisan@m.com"
$ws.Range("D13").Value = "+8613512245671"
$ws.Range("E13").Value = "公司/部门C/中心CA/小组CAA"
$ws.Range("F13").Value = "lisi"
$ws.Range("G13").Value = 30
$ws.Range("H13").Value = "female"
$ws.Range("I13").Value = "region-10"

# 4. Update the "freedom" row (now row 14) with its new age/gender/region values
$ws.Range("G14").Value = 666
$ws.Range("H14").Value = "other"
$ws.Range("I14").Value = "solar-system"

# 5. Apply the text number format across the whole data range (matches style index 1 "@")
$ws.Range("A3:I14").NumberFormat = "@"

# 6. Turn the new C13 email into a mailto hyperlink
$ws.Hyperlinks.Add($ws.Range("C13"), "mailto:qinshisan@m.com", "", "", "qinshisan@m.com")

# 7. Update the active selection to reflect where the author ended up
$ws.Range("H16").Select()
